$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add a new "2022-Q3" sheet, positioned right before "2022-Q2",
#    by copying the "2022-Q2" sheet (so headers/formatting match the
#    other quarterly detail sheets) and overwriting its data row.
# ---------------------------------------------------------------------
$srcQ2 = $wb.Worksheets.Item("2022-Q2")
$srcQ2.Copy($srcQ2)
$newSheet = $wb.ActiveSheet
$newSheet.Name = "2022-Q3"

# Force the numeric-looking data cells to stay text (matches source
# sheet's inline string cells - otherwise strings like "159726" or
# "2.60" get auto-coerced to numbers and lose formatting/leading data).
# (C2, the fund name, is unambiguous text already and needs no help.)
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("D2:G2").NumberFormat = "@"

$newSheet.Range("B2").Value = "159726"
$newSheet.Range("C2").Value = "华夏恒生中国内地企业高股息率ETF"
$newSheet.Range("D2").Value = "0.84"
$newSheet.Range("E2").Value = "96.48"
$newSheet.Range("F2").Value = "2.60"
$newSheet.Range("G2").Value = "0.0218"
$newSheet.Range("H2").Value = 9

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q3
#    right under the header, push the existing quarters down, and
#    keep the running index in column A sequential.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.02

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# Re-apply the plain data-row formatting (the row-insert above copies
# the preceding row's format onto the new row, which picks up bold /
# border styling that the other data rows don't have).
$summary.Range("A3:D3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# Restore the original tab selection (the "2020-Q4" sheet was the
# active/selected tab before this edit).
$wb.Worksheets.Item("2020-Q4").Activate()
